$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "token[code]"
$ws.Range("B17").Value = 32
$ws.Range("C17").Value = 32
$ws.Range("F17").Value = "^[0-9a-f]{32}$"

$ws.Range("F17").Select()
